$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = [double]"0.008496645711067738"
$ws.Range("B3").Value = [double]"0.772939363834608"
$ws.Range("B4").Value = [double]"0.1048964807954936"
$ws.Range("B5").Value = "[1.0, 0.020993214173052173, -0.04176833174737117, -0.10560999193303837, -0.11533516583456446, 0.07035599977997269, 0.026391251168427818, 0.039913619224905376, 0.020689207515053064, -0.01731142237911526, -0.051325861500777385, -0.03442789447122897, 0.0011335207199031814, 0.059519924174354234, 0.056167899344512606, 0.032211846857381946, -0.01283971000304848, -0.031036381260136484, -0.03772868825354034, -0.002959421961216216]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = [double]"0.8847643641058748"
$ws.Range("B4").Value = [double]"-0.1812764579519923"
$ws.Range("B5").Value = [double]"0.5458046833958783"
$ws.Range("B6").Value = [double]"0.06305607484792758"
$ws.Range("B7").Value = [double]"0.1799868863271274"
$ws.Range("B8").Value = [double]"0.5749375572395075"
$ws.Range("B9").Value = [double]"0.06303061886141828"
$ws.Range("B10").Value = "[1.0, 0.05936828066251595, 0.01378317064568642, 0.0013389214288225123, -0.00704774752881248, 0.0019225910247095503, 0.010772661747817163, -0.023827857189336794, 0.0032501902115645103, -0.012778729844369953, 0.004055392072017083, 0.014332596877540838, 0.012405963046405833, 0.03059605823726628, 0.029353599051395, 0.0036666200381629196, -0.01638440717766988, -0.005631955863527286, -0.0010197583031796833, -0.011009909056893315]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = [double]"0.0009873166335556474"
$ws.Range("B3").Value = [double]"0.1133115976320104"
$ws.Range("B4").Value = [double]"0.04222293798008307"
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.0035827774405508667, -0.03340950229404528, -0.09611604744631938, -0.12828549937149172, 0.09389995351401186, 0.047924417783650544, 0.040231768982029424, 0.023845263578730742, -0.009572706618660347, -0.057195874014094296, -0.028917465885140933, 0.0006093134134932231, 0.05672645911762163, 0.04550471288916477, 0.019594759359999013, -0.00822015815457666, -0.028157406905283143, -0.05026168883274082, -0.0043637928357312394]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = [double]"-0.0008399168189114137"
$ws.Range("B3").Value = [double]"0.113002184732998"
$ws.Range("B4").Value = [double]"0.01322364626117355"
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.002532115242903328, -0.03374763094583777, -0.09591283751275625, -0.12865656491069777, 0.0943613200162274, 0.04798658597518958, 0.03923494576515497, 0.02460823006793059, -0.010023409121428005, -0.05636285875729877, -0.029063338915024616, 0.0005479569371484508, 0.05645980672459899, 0.04521190515642551, 0.019553193672920893, -0.008112357540072222, -0.02856324722072384, -0.050440185404500223, -0.004607816837164434]"
$ws.Range("B7").Value = [double]"0.06340350266431082"

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = [double]"0.008576308576116158"
$ws.Range("B3").Value = [double]"0.1024876884301354"
$ws.Range("B4").Value = [double]"1.805300617031429e-10"
$ws.Range("B6").Value = "[0.9999999999999998, 0.0008215369934537743, -0.04033836938288315, -0.10505544738136322, -0.11536158852912097, 0.07311314602788013, 0.024833604709608496, 0.03913659028049754, 0.02019155537791576, -0.01647136314383136, -0.049363228480920236, -0.03398459019177501, 0.0007680003207134622, 0.057600756789735724, 0.054690747790072644, 0.03141910170917143, -0.012750466758428274, -0.030560715058727473, -0.03727699042133201, -0.0027492667194408847]"
$ws.Range("B7").Value = [double]"0.04526849715481648"
$ws.Range("B9").Value = [double]"0.77383177592694"
